# Rename the report template placeholder used in the "expanded expense
# analysis" table header (cell B10) from
#   {$v->rows[]->expense_label}
# to
#   {$v->rows[]->label}
# Single-quoted string so PowerShell does not try to expand "$v" as a
# variable reference.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B10").Value = '{$v->rows[]->label}'
